$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 data: A7 = "a1"
$ws.Cells.Item(7, 1).Value = "a1"
$ws.Cells.Item(7, 2).Value = 0.8709480166435242
$ws.Cells.Item(7, 3).Value = 0.8732147216796875
$ws.Cells.Item(7, 4).Value = 269.6632690429688
$ws.Cells.Item(7, 5).Value = 32.62083053588867
$ws.Cells.Item(7, 6).Value = 32.26833343505859
$ws.Cells.Item(7, 7).Value = 181.3432159423828
$ws.Cells.Item(7, 8).Value = 228.9816131591797

# New row 8 data: A8 = "b2"
$ws.Cells.Item(8, 1).Value = "b2"
$ws.Cells.Item(8, 2).Value = 0.8709480166435242
$ws.Cells.Item(8, 3).Value = 0.8732147216796875
$ws.Cells.Item(8, 4).Value = 269.6569213867188
$ws.Cells.Item(8, 5).Value = 32.6202278137207
$ws.Cells.Item(8, 6).Value = 32.26775360107422
$ws.Cells.Item(8, 7).Value = 181.3432006835938
$ws.Cells.Item(8, 8).Value = 228.9815063476562

# Copy the style used by the other A-column label cells (e.g. A6) to A7 and A8
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
